# Update "想去人数" (F column) figures across the relevant worksheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8894
$ws1.Range("F3").Value = 97
$ws1.Range("F5").Value = 103
$ws1.Range("F6").Value = 1483
$ws1.Range("F7").Value = 1406
$ws1.Range("F8").Value = 248
$ws1.Range("F10").Value = 318
$ws1.Range("F11").Value = 90

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8894
$ws4.Range("F3").Value = 97
$ws4.Range("F5").Value = 103
$ws4.Range("F6").Value = 1483
$ws4.Range("F7").Value = 1406
$ws4.Range("F8").Value = 248
$ws4.Range("F11").Value = 318
$ws4.Range("F12").Value = 90
